$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 43
$ws.Range("H43").Value = 1474.4706
$ws.Range("I43").Value = 1031.3334
$ws.Range("J43").Value = 1716.1818
$ws.Range("K43").Value = 1031.3334
$ws.Range("L43").Value = 1716.1818
$ws.Range("M43").Value = -962.3334
$ws.Range("N43").Value = -1854.1818

# ALC row 64
$ws.Range("H64").Value = 3846.6667
$ws.Range("I64").Value = 4122.857
$ws.Range("J64").Value = 3460
$ws.Range("K64").Value = 4122.857
$ws.Range("L64").Value = 3460
$ws.Range("M64").Value = -3874.857
$ws.Range("N64").Value = -3956

# ALC row 67
$ws.Range("H67").Value = 3846.6667
$ws.Range("I67").Value = 4122.857
$ws.Range("J67").Value = 3460
$ws.Range("K67").Value = 4122.857
$ws.Range("L67").Value = 3460
$ws.Range("M67").Value = -3264.857
$ws.Range("N67").Value = -5176

# ALC row 76
$ws.Range("H76").Value = 3100
$ws.Range("I76").Value = 2700
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 2700
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -2385
$ws.Range("N76").Value = -5630

# ALC row 79
$ws.Range("H79").Value = 3100
$ws.Range("I79").Value = 2700
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 2700
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -1608
$ws.Range("N79").Value = -7184

$ws = $wb.Worksheets.Item("ARM")
# ARM row 63
$ws.Range("H63").Value = 2371.0527
$ws.Range("I63").Value = 2330.9858
$ws.Range("J63").Value = 2940
$ws.Range("K63").Value = 2330.9858
$ws.Range("L63").Value = 2940
$ws.Range("M63").Value = -1644.9858
$ws.Range("N63").Value = -4312

# ARM row 66
$ws.Range("H66").Value = 2371.0527
$ws.Range("I66").Value = 2330.9858
$ws.Range("J66").Value = 2940
$ws.Range("K66").Value = 11654.929
$ws.Range("L66").Value = 14700
$ws.Range("M66").Value = -8222.929
$ws.Range("N66").Value = -21564

# ARM row 88
$ws.Range("H88").Value = 2906
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2906
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2906
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -3718

# ARM row 91
$ws.Range("H91").Value = 2906
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2906
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2906
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -5714

# ARM row 122
$ws.Range("H122").Value = 8930221
$ws.Range("I122").Value = 13890322
$ws.Range("J122").Value = 2040
$ws.Range("K122").Value = 41670966
$ws.Range("L122").Value = 6120
$ws.Range("M122").Value = -41668516
$ws.Range("N122").Value = -11020

# ARM row 132
$ws.Range("H132").Value = 1693.0377
$ws.Range("I132").Value = 1212.3
$ws.Range("J132").Value = 3172.2307
$ws.Range("K132").Value = 3636.9
$ws.Range("L132").Value = 9516.6921
$ws.Range("M132").Value = -1106.9
$ws.Range("N132").Value = -14576.6921

$ws = $wb.Worksheets.Item("BSM")
# BSM row 105
$ws.Range("H105").Value = 1625521.5
$ws.Range("I105").Value = 3248314.5
$ws.Range("J105").Value = 2728.5715
$ws.Range("K105").Value = 3248314.5
$ws.Range("L105").Value = 2728.5715
$ws.Range("M105").Value = -3246567.5
$ws.Range("N105").Value = -6222.5715

# BSM row 133
$ws.Range("H133").Value = 35534.832
$ws.Range("I133").Value = 70709
$ws.Range("J133").Value = 28500
$ws.Range("K133").Value = 70709
$ws.Range("L133").Value = 28500
$ws.Range("M133").Value = -65649
$ws.Range("N133").Value = -38620

$ws = $wb.Worksheets.Item("CRP")
# CRP row 62
$ws.Range("H62").Value = 2989.1765
$ws.Range("I62").Value = 2676.25
$ws.Range("J62").Value = 3267.3333
$ws.Range("K62").Value = 2676.25
$ws.Range("L62").Value = 3267.3333
$ws.Range("M62").Value = -2052.25
$ws.Range("N62").Value = -4515.3333

# CRP row 65
$ws.Range("H65").Value = 2989.1765
$ws.Range("I65").Value = 2676.25
$ws.Range("J65").Value = 3267.3333
$ws.Range("K65").Value = 13381.25
$ws.Range("L65").Value = 16336.6665
$ws.Range("M65").Value = -10261.25
$ws.Range("N65").Value = -22576.6665

# CRP row 122
$ws.Range("H122").Value = 2777.8462
$ws.Range("I122").Value = 2003
$ws.Range("J122").Value = 3122.2222
$ws.Range("K122").Value = 6009
$ws.Range("L122").Value = 9366.6666
$ws.Range("M122").Value = -3559
$ws.Range("N122").Value = -14266.6666

# CRP row 132
$ws.Range("H132").Value = 1954.9535
$ws.Range("I132").Value = 1474.5
$ws.Range("J132").Value = 2458.2856
$ws.Range("K132").Value = 4423.5
$ws.Range("L132").Value = 7374.8568
$ws.Range("M132").Value = -1893.5
$ws.Range("N132").Value = -12434.8568

$ws = $wb.Worksheets.Item("CUL")
# CUL row 86
$ws.Range("H86").Value = 512.2778
$ws.Range("I86").Value = 490.16666
$ws.Range("J86").Value = 556.5
$ws.Range("K86").Value = 1470.49998
$ws.Range("L86").Value = 1669.5
$ws.Range("M86").Value = -284.4999800000001
$ws.Range("N86").Value = -4041.5

# CUL row 89
$ws.Range("H89").Value = 512.2778
$ws.Range("I89").Value = 490.16666
$ws.Range("J89").Value = 556.5
$ws.Range("K89").Value = 4411.49994
$ws.Range("L89").Value = 5008.5
$ws.Range("M89").Value = 1516.50006
$ws.Range("N89").Value = -16864.5

# CUL row 131
$ws.Range("H131").Value = 11364752
$ws.Range("I131").Value = 83333590
$ws.Range("J131").Value = 1251.3158
$ws.Range("K131").Value = 250000770
$ws.Range("L131").Value = 3753.9474
$ws.Range("M131").Value = -249995730
$ws.Range("N131").Value = -13833.9474

# CUL row 136
$ws.Range("H136").Value = 800
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 800
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 2400
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("GSM")
# GSM row 70
$ws.Range("H70").Value = 4306.9546
$ws.Range("I70").Value = 4225.905
$ws.Range("J70").Value = 6009
$ws.Range("K70").Value = 4225.905
$ws.Range("L70").Value = 6009
$ws.Range("M70").Value = -3955.905
$ws.Range("N70").Value = -6549

# GSM row 73
$ws.Range("H73").Value = 4306.9546
$ws.Range("I73").Value = 4225.905
$ws.Range("J73").Value = 6009
$ws.Range("K73").Value = 4225.905
$ws.Range("L73").Value = 6009
$ws.Range("M73").Value = -3289.905
$ws.Range("N73").Value = -7881

# GSM row 80
$ws.Range("H80").Value = 3899.5454
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3899.5454
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3899.5454
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -5895.5454

# GSM row 83
$ws.Range("H83").Value = 3899.5454
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3899.5454
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 19497.727
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -29481.727

$ws = $wb.Worksheets.Item("LTW")
# LTW row 16
$ws.Range("H16").Value = 1909.3846
$ws.Range("I16").Value = 1943.5
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1943.5
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -1773.5
$ws.Range("N16").Value = -1840

# LTW row 40
$ws.Range("H40").Value = 2586.6875
$ws.Range("I40").Value = 1833.5555
$ws.Range("J40").Value = 3555
$ws.Range("K40").Value = 1833.5555
$ws.Range("L40").Value = 3555
$ws.Range("M40").Value = -1697.5555
$ws.Range("N40").Value = -3827

$ws = $wb.Worksheets.Item("WVR")
# WVR row 61
$ws.Range("H61").Value = 9003
$ws.Range("I61").Value = 1525.5
$ws.Range("J61").Value = 10872.375
$ws.Range("K61").Value = 1525.5
$ws.Range("L61").Value = 10872.375
$ws.Range("M61").Value = -1233.5
$ws.Range("N61").Value = -11456.375
